$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Data rows (A:E) for rows 2..12 - new set of sample Excel/Json file paths.
# Column A = source .xlsx path, Column B = output .json path,
# Column C = header_row, Column D = data_row, Column E = start_col.
# ---------------------------------------------------------------------------
$rows = @(
    @("D:\PythonProjects\Jexcel\Samples\Excels1\ChainedList.1.xlsx",             "D:\PythonProjects\Jexcel\Samples\OutJson1\ChainedList.1.json"),
    @("D:\PythonProjects\Jexcel\Samples\Excels1\ChainedList.1-var.1.xlsx",       "D:\PythonProjects\Jexcel\Samples\OutJson1\ChainedList.1-var.1.json"),
    @("D:\PythonProjects\Jexcel\Samples\Excels1\Simple_Dict_List.1.xlsx",        "D:\PythonProjects\Jexcel\Samples\OutJson1\Simple_Dict_List.1.json"),
    @("D:\PythonProjects\Jexcel\Samples\Excels1\Simple_Dict_List.1-Complex.xlsx","D:\PythonProjects\Jexcel\Samples\OutJson1\Simple_Dict_List.1-Complex.json"),
    @("D:\PythonProjects\Jexcel\Samples\Excels1\Simple_Dict_List.2.xlsx",        "D:\PythonProjects\Jexcel\Samples\OutJson1\Simple_Dict_List.2.json"),
    @("D:\PythonProjects\Jexcel\Samples\Excels1\Simple_Dict_List.3.xlsx",        "D:\PythonProjects\Jexcel\Samples\OutJson1\Simple_Dict_List.3.json"),
    @("D:\PythonProjects\Jexcel\Samples\Excels1\Simple_Dict_List.4.xlsx",        "D:\PythonProjects\Jexcel\Samples\OutJson1\Simple_Dict_List.4.json"),
    @("D:\PythonProjects\Jexcel\Samples\Excels1\Simple_Dict_List.5.xlsx",        "D:\PythonProjects\Jexcel\Samples\OutJson1\Simple_Dict_List.5.json"),
    @("D:\PythonProjects\Jexcel\Samples\Excels1\Simple_Dict_List.5-var.1.xlsx",  "D:\PythonProjects\Jexcel\Samples\OutJson1\Simple_Dict_List.5-var.1.json"),
    @("D:\PythonProjects\Jexcel\Samples\Excels1\Simple_Dict_List.5-var.2.xlsx",  "D:\PythonProjects\Jexcel\Samples\OutJson1\Simple_Dict_List.5-var.2.json"),
    @("D:\PythonProjects\Jexcel\Samples\Excels2\1.xlsx",                        "D:\PythonProjects\Jexcel\Samples\OutJson2\1.json")
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = 0
    $ws.Cells.Item($r, 4).Value = -1
    $ws.Cells.Item($r, 5).Value = 0
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Column A/B of the new rows get a "left, centered, indent 4" look (matches
# the new cellXfs entry). Build the style cleanly on A2 (one of the cells we
# just wrote) and format-paint it across A2:B12 so no stray intermediate
# styles get minted.
# ---------------------------------------------------------------------------
$seed = $ws.Range("A2")
$seed.HorizontalAlignment = -4131   # xlLeft
$seed.VerticalAlignment = -4108     # xlCenter
$seed.IndentLevel = 4
$seed.Copy()
$ws.Range("A2:B12").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# Columns C:E for the newly added rows 5..12 keep the original
# centered-style (already used by row 1 / rows 2-4) - paint it from C1.
$seedCE = $ws.Range("C1")
$seedCE.Copy()
$ws.Range("C2:E12").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Column widths (characters) - A & B grew wider to fit the longer paths.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 70.5
$ws.Columns.Item(2).ColumnWidth = 69.375

# ---------------------------------------------------------------------------
# Selection moved to B17 (author scrolled/selected past the new data).
# ---------------------------------------------------------------------------
$ws.Range("B17").Select()

# ---------------------------------------------------------------------------
# Window was maximized (xWindow/yWindow -120, larger width/height).
# ---------------------------------------------------------------------------
$excel.ActiveWindow.WindowState = -4137   # xlMaximized
